$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the student username/email columns from the old IIIT Allahabad
# naming scheme to the new Thapar scheme.
$ws.Range("B2").Value = "IT2050001"
$ws.Range("C2").Value = "it2050001@thapar.edu"

$ws.Range("B3").Value = "IT2050002"
$ws.Range("C3").Value = "it2050002@thapar.edu"

$ws.Range("B4").Value = "IT2050003"
$ws.Range("C4").Value = "it2050003@thapar.edu"

$ws.Range("B5").Value = "IT2050004"
$ws.Range("C5").Value = "it2050004@thapar.edu"

$ws.Range("B6").Value = "IT2050005"
$ws.Range("C6").Value = "it2050005@thapar.edu"

# Match the last active selection recorded in the saved file.
$ws.Range("C7").Select()
